$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.0
$ws.Range("H2").Value = 0.1025354202104953
$ws.Range("B3").Value = 0.09438849744921556
$ws.Range("H3").Value = 0.1969239176597109
$ws.Range("B4").Value = 0.07262270252128764
$ws.Range("H4").Value = 0.1751581227317829
$ws.Range("B5").Value = 0.04880804554619359
$ws.Range("H5").Value = 0.1513434657566889
$ws.Range("B6").Value = 0.03060681022136408
$ws.Range("H6").Value = 0.1331422304318594
$ws.Range("B7").Value = 0.01967513280245653
$ws.Range("C7").Value = 0.002401627296344176
$ws.Range("D7").Value = 2.025554899249398
$ws.Range("E7").Value = 0.01094939379502432
$ws.Range("F7").Value = 0.01495807827251136
$ws.Range("G7").Value = 0.02439218733240197
$ws.Range("H7").Value = 0.1222105530129518
$ws.Range("B8").Value = 0.01880783187435991
$ws.Range("C8").Value = 0.002530411786938378
$ws.Range("D8").Value = 2.015292309059852
$ws.Range("E8").Value = 0.007783556192342187
$ws.Range("F8").Value = 0.01382927390775375
$ws.Range("G8").Value = 0.02378638984096626
$ws.Range("H8").Value = 0.1213432520848552
$ws.Range("B9").Value = 0.01782325126250029
$ws.Range("C9").Value = 0.002140609952760626
$ws.Range("D9").Value = 1.774552441494788
$ws.Range("E9").Value = 0.004373944970796132
$ws.Range("F9").Value = 0.01362314459485887
$ws.Range("G9").Value = 0.02202335793014191
$ws.Range("H9").Value = 0.1203586714729956
$ws.Range("B10").Value = 0.01781243772760065
$ws.Range("C10").Value = 0.001993052578255619
$ws.Range("D10").Value = 1.608405485260382
$ws.Range("E10").Value = 0.005759683438276856
$ws.Range("F10").Value = 0.01389698689346763
$ws.Range("G10").Value = 0.02172788856173453
$ws.Range("H10").Value = 0.120347857938096
$ws.Range("B11").Value = 0.03045470186181892
$ws.Range("H11").Value = 0.1329901220723142
$ws.Range("B12").Value = 0.05251030857758803
$ws.Range("H12").Value = 0.1550457287880833
$ws.Range("B13").Value = 0.06517608269330417
$ws.Range("H13").Value = 0.1677115029037995
$ws.Range("B14").Value = 0.07306994158361906
$ws.Range("H14").Value = 0.1756053617941144
$ws.Range("B15").Value = 0.08049854346949885
$ws.Range("H15").Value = 0.1830339636799942
$ws.Range("B16").Value = 0.08373282897400397
$ws.Range("H16").Value = 0.1862682491844993
$ws.Range("B17").Value = 0.08495449961570588
$ws.Range("H17").Value = 0.1874899198262012
$ws.Range("B18").Value = -0.1025354202104953
$ws.Range("C18").Value = 0.008773877970870389
$ws.Range("D18").Value = -147864385954.2316
$ws.Range("E18").Value = 0.03158185907953013
$ws.Range("F18").Value = -0.1197770267284051
$ws.Range("G18").Value = -0.08529381369258568
$ws.Range("H18").Value = 0.0
$ws.Range("B19").Value = 0.08524194045651438
$ws.Range("H19").Value = 0.1877773606670097
$ws.Range("B20").Value = 0.09045507513899632
$ws.Range("H20").Value = 0.1929904953494916
$ws.Range("B21").Value = 0.09267026568961423
$ws.Range("H21").Value = 0.1952056859001096
$ws.Range("B22").Value = 0.0944551338105561
$ws.Range("H22").Value = 0.1969905540210514
$ws.Range("B23").Value = 0.1003439645995695
$ws.Range("C23").Value = 0.007970831847943354
$ws.Range("D23").Value = 1596169204495.518
$ws.Range("E23").Value = 0.04506884183514018
$ws.Range("F23").Value = 0.08466526849088823
$ws.Range("G23").Value = 0.116022660708251
$ws.Range("H23").Value = 0.2028793848100648
$ws.Range("B24").Value = 0.1016688522498076
$ws.Range("C24").Value = 0.00761966757086134
$ws.Range("D24").Value = 970354239318.228
$ws.Range("E24").Value = 0.04728213471158763
$ws.Range("F24").Value = 0.08669501588202619
$ws.Range("G24").Value = 0.1166426886175892
$ws.Range("H24").Value = 0.2042042724603029
$ws.Range("B25").Value = 0.1030306421282219
$ws.Range("C25").Value = 0.00810560903366833
$ws.Range("D25").Value = 3111251340495.385
$ws.Range("E25").Value = 0.05795883556857434
$ws.Range("F25").Value = 0.08708429616865478
$ws.Range("G25").Value = 0.1189769880877886
$ws.Range("H25").Value = 0.2055660623387172
$ws.Range("B26").Value = 0.1052088270725087
$ws.Range("H26").Value = 0.207744247283004
$ws.Range("B27").Value = 0.1056437687550505
$ws.Range("C27").Value = 0.00829279257810839
$ws.Range("D27").Value = 1126179147947.119
$ws.Range("E27").Value = 0.07498685944630981
$ws.Range("F27").Value = 0.08933459140845156
$ws.Range("G27").Value = 0.1219529461016493
$ws.Range("H27").Value = 0.2081791889655458
$ws.Range("B28").Value = 0.102342716796393
$ws.Range("C28").Value = 0.007681622173093967
$ws.Range("D28").Value = 20.44270115233925
$ws.Range("E28").Value = 0.1177495005018594
$ws.Range("F28").Value = 0.08725755755123202
$ws.Range("G28").Value = 0.1174278760415537
$ws.Range("H28").Value = 0.2048781370068883
$ws.Range("B29").Value = 0.01956145910228603
$ws.Range("C29").Value = 0.002134752326741195
$ws.Range("D29").Value = 2.201256146179134
$ws.Range("E29").Value = 0.009606414221899074
$ws.Range("F29").Value = 0.01535772402258803
$ws.Range("G29").Value = 0.02376519418198398
$ws.Range("H29").Value = 0.1220968793127813
